$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (llama-3.1-8b-instant)
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 0.001111111111111111
$ws.Range("K2").Value = 4868
$ws.Range("L2").Value = 0.009736

# Row 3 (llama-3.3-70b-versatile)
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 0.002
$ws.Range("K3").Value = 712
$ws.Range("L3").Value = 0.00712
